$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the username list (A2:A5) with the corrected values.
$ws.Range("A2").Value = "__sanatani__090"
$ws.Range("A3").Value = "_abo_safwan_"
$ws.Range("A4").Value = "__r_n_shanawar__"
$ws.Range("A5").Value = "_agencia.dara"

# Update the selected cell to match the saved view state.
$ws.Range("H4").Select()
